# Daily attendance processing - 2025-10-25 06:26:35
#
# Reorders the "Recorded By" (column G) value on the "Session Analysis
# Results" sheet: any "System" / "system" token(s) currently at the front
# of the comma-separated recorder list are moved to the end, so the real
# user/email recorders sort first (e.g. "System, dnasr281@gmail.com"
# becomes "dnasr281@gmail.com, System"). Rows whose recorder list still
# includes "admin@admin.com" are left untouched, matching upstream data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1
$col = 7   # column G = "Recorded By"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }
    if ($val -notmatch ",") {
        continue
    }
    if ($val -match "admin@admin.com") {
        continue
    }

    $parts = $val -split ", "

    $systemParts = @()
    $otherParts = @()
    foreach ($part in $parts) {
        if ($part.ToLower() -eq "system") {
            $systemParts += $part
        } else {
            $otherParts += $part
        }
    }

    if ($systemParts.Count -eq 0) {
        continue
    }

    $newParts = $otherParts + $systemParts
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
